$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1 should carry the same formatting (bold font,
# border, centered/top alignment) as the existing header cells (e.g. H1).
# Copy H1's format onto I1:J1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J (plain numeric cells, no special style).
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 6
